$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(4.0204,4.5157,1.6098),
    @(3.2872,3.8755,1.8316),
    @(3.4021,3.9344,1.5457),
    @(3.0174,3.382,1.5601),
    @(3.045,3.6609,1.8298),
    @(2.8501,3.3194,1.5275),
    @(3.1895,3.4672,1.5194),
    @(2.9765,3.4629,1.8444),
    @(2.7317,3.2242,1.4422),
    @(3.3949,3.5549,1.319),
    @(3.1364,3.6534,2.2688),
    @(2.4364,2.8507,1.2094),
    @(2.6353,2.7999,0.8816000000000001),
    @(2.659,2.9465,1.0792),
    @(2.2528,2.4957,0.8583),
    @(2.2528,2.4957,0.8583),
    @(2.659,2.9465,1.0792),
    @(2.6353,2.7999,0.8816000000000001),
    @(2.4364,2.8507,1.2094),
    @(3.1364,3.6534,2.2688),
    @(3.3949,3.5549,1.319),
    @(2.7317,3.2242,1.4422),
    @(2.9765,3.4629,1.8444),
    @(3.1895,3.4672,1.5194),
    @(2.8501,3.3194,1.5275),
    @(3.045,3.6609,1.8298),
    @(3.0174,3.382,1.5601),
    @(3.4021,3.9344,1.5457),
    @(3.2872,3.8755,1.8316),
    @(4.0204,4.5157,1.6098),
    @(4.4896,4.1077,15.7036),
    @(4.939,3.9113,14.8169),
    @(4.3764,4.0127,11.5296),
    @(4.2893,4.5044,17.4173),
    @(4.0908,4.2797,12.2723),
    @(4.0948,4.2905,9.832599999999999),
    @(4.2827,4.5327,17.2035),
    @(4.1559,4.4095,12.4676),
    @(4.0943,4.3221,9.8894),
    @(5.0227,5.3436,12.8589),
    @(4.1972,4.4192,12.5915),
    @(3.6008,3.7628,4.2247),
    @(4.2025,4.4195,12.6076),
    @(4.2019,4.4162,12.6058),
    @(4.2025,4.4195,12.6076),
    @(3.6008,3.7628,4.2247),
    @(4.1972,4.4192,12.5915),
    @(5.0227,5.3436,12.8589),
    @(4.0943,4.3221,9.8894),
    @(4.1559,4.4095,12.4676),
    @(4.2827,4.5327,17.2035),
    @(4.0948,4.2905,9.832599999999999),
    @(4.0908,4.2797,12.2723),
    @(4.2893,4.5044,17.4173),
    @(4.3764,4.0127,11.5296),
    @(4.939,3.9113,14.8169),
    @(4.4896,4.1077,15.7036)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $r = $startRow + $i
    $row = $newValues[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
}
